# Refresh cryptos.xlsx price/volume data (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.167.60'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.60%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.787.26'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.79%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '619.83'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.75'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.43%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.784.05'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.87%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.536'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.99%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.26%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.34'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.73%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.493'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.22'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.57%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000263'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.46%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.425.28'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.790.66'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.209.73'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.66%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.65'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.40%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.84'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.87%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '511.82'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.66'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.729'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.52'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.92%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '88.16'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.22'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.77%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.10'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000140'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +24.95%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.50'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.58%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.86'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.13%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.77'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -6.28%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.57'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.15%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.24%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.19'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.34%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.336'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.93%  '

$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.133'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.22%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.15'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.08'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '44.57'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.86%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.52%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '423.14'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.058.22'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.77'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.77%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0364'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.70'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.15%  '

$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.02%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '137.02'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.45%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.06%  '
